$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rewrite "ODI Batting" sheet (currently sheet index 1):
#    - header D1: MATCH_CARD_LINK -> MATCH_CODE
#    - column D values: full howstat URL -> trailing MatchCode number (as text)
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item(1)

$battingSheet.Range("D1").Value = "MATCH_CODE"

for ($r = 2; $r -le 219; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $txt = $cell.Text
    $eqIdx = $txt.LastIndexOf("=")
    if ($eqIdx -ge 0) {
        $code = $txt.Substring($eqIdx + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $code
    }
}

# ---------------------------------------------------------------------------
# 2. Rewrite "ODI Bowling" sheet (currently sheet index 2):
#    - header B1: MATCH_CARD_LINK -> MATCH_CODE
#    - column B values: full howstat URL -> trailing MatchCode number (as text)
# ---------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item(2)

$bowlingSheet.Range("B1").Value = "MATCH_CODE"

for ($r = 2; $r -le 147; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $txt = $cell.Text
    $eqIdx = $txt.LastIndexOf("=")
    if ($eqIdx -ge 0) {
        $code = $txt.Substring($eqIdx + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $code
    }
}

# ---------------------------------------------------------------------------
# 3. Insert a new "Player Info" sheet as the very first sheet.
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 0; $c -lt $piHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $piHeaders[$c]
}

$piRow = @("3488", "Mohammad Mahmudullah", "Right Handed", "Right Arm Off Break")
for ($c = 0; $c -lt $piRow.Length; $c++) {
    $cell = $playerInfo.Cells.Item(2, $c + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $piRow[$c]
}

# ---------------------------------------------------------------------------
# 4. Append a new "ODI Batting Extra" sheet as the very last sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add($null, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 0; $c -lt $extraHeaders.Length; $c++) {
    $cell = $battingExtra.Cells.Item(1, $c + 1)
    $cell.NumberFormat = "@"
    $cell.Value = $extraHeaders[$c]
}

# Each inner element is @(value, type) where type is "n" (real number) or "s" (text).
$extraData = @(
    ,@(@("4479","s"), @("6","n"), @("3","s"), @("0","s"), @("10.74%","s"), @("NO","s"))
    ,@(@("4481","s"), @("5","n"), @("0","s"), @("0","s"), @("","s"), @("NO","s"))
    ,@(@("4537","s"), @("6","n"), @("1","s"), @("0","s"), @("3.65%","s"), @("NO","s"))
    ,@(@("4538","s"), @("","s"), @("","s"), @("","s"), @("","s"), @("NO","s"))
    ,@(@("4539","s"), @("6","n"), @("0","s"), @("0","s"), @("15.10%","s"), @("NO","s"))
    ,@(@("4550","s"), @("6","n"), @("1","s"), @("1","s"), @("7.96%","s"), @("NO","s"))
    ,@(@("4557","s"), @("6","n"), @("3","s"), @("0","s"), @("12.89%","s"), @("NO","s"))
    ,@(@("4559","s"), @("6","n"), @("","s"), @("","s"), @("","s"), @("NO","s"))
    ,@(@("4606","s"), @("","s"), @("","s"), @("","s"), @("","s"), @("NO","s"))
    ,@(@("4611","s"), @("4","n"), @("","s"), @("","s"), @("","s"), @("NO","s"))
    ,@(@("4616","s"), @("4","n"), @("1","s"), @("0","s"), @("13.41%","s"), @("NO","s"))
    ,@(@("4626","s"), @("5","n"), @("3","s"), @("0","s"), @("6.60%","s"), @("NO","s"))
    ,@(@("4627","s"), @("5","n"), @("3","s"), @("3","s"), @("27.59%","s"), @("NO","s"))
    ,@(@("4628","s"), @("5","n"), @("3","s"), @("0","s"), @("15.23%","s"), @("NO","s"))
    ,@(@("4679","s"), @("6","n"), @("0","s"), @("0","s"), @("7.49%","s"), @("NO","s"))
    ,@(@("4682","s"), @("6","n"), @("7","s"), @("0","s"), @("28.41%","s"), @("NO","s"))
    ,@(@("4685","s"), @("","s"), @("","s"), @("","s"), @("","s"), @("NO","s"))
    ,@(@("4711","s"), @("6","n"), @("3","s"), @("0","s"), @("14.83%","s"), @("NO","s"))
    ,@(@("4713","s"), @("6","n"), @("3","s"), @("0","s"), @("16.49%","s"), @("NO","s"))
    ,@(@("4717","s"), @("6","n"), @("0","s"), @("1","s"), @("3.25%","s"), @("NO","s"))
)

for ($i = 0; $i -lt $extraData.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $extraData[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $pair = $rowVals[$c]
        $val = $pair[0]
        $kind = $pair[1]
        $cell = $battingExtra.Cells.Item($rowNum, $c + 1)
        if ($kind -eq "n") {
            $cell.Value = [double]$val
        } else {
            $cell.NumberFormat = "@"
            $cell.Value = $val
        }
    }
}

Write-Host "edit complete"
